# Update the "想去人数" (F column) values across the sheets, per the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 2-23
$ws1 = $wb.Worksheets.Item("展览")
$values1 = @(117,208,6,6551,78,428,131,5911,38,190,1234,7,82,382,87,15,339,39,3,4228,35,182)
for ($i = 0; $i -lt $values1.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $values1[$i]
}

# Sheet "演出" (sheet2): row 2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 43

# Sheet "本地生活" (sheet3): no data rows to update (only header row present)

# Sheet "全部类型" (sheet4): rows 2-24
$ws4 = $wb.Worksheets.Item("全部类型")
$values4 = @(117,208,6,6551,78,428,131,5911,38,190,1234,7,82,382,87,15,339,39,3,4228,43,35,182)
for ($i = 0; $i -lt $values4.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 6).Value = $values4[$i]
}
